$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 12.92888040645593
$ws.Range("C2").Value = 9.070380523056762
$ws.Range("E2").Value = 11.43361080149588
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 28.70288118982866
$ws.Range("H2").Value = 13.68321391459183
$ws.Range("M2").Value = 14.64663356233735
$ws.Range("B3").Value = 12.31466191686637
$ws.Range("C3").Value = 8.596459605942661
$ws.Range("E3").Value = 11.3357051710136
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 28.40597409033242
$ws.Range("H3").Value = 13.72347546076177
$ws.Range("M3").Value = 14.35755425966472
$ws.Range("B4").Value = 11.92375904867162
$ws.Range("C4").Value = 8.290032952169081
$ws.Range("E4").Value = 11.27936703737528
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 28.24018911060421
$ws.Range("H4").Value = 13.75268308668547
$ws.Range("M4").Value = 14.18108649486832
$ws.Range("B5").Value = 11.76122278339534
$ws.Range("C5").Value = 8.161344771634434
$ws.Range("E5").Value = 11.25737832184468
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 28.17685322940528
$ws.Range("H5").Value = 13.76570531704673
$ws.Range("M5").Value = 14.10953694384014
$ws.Range("B6").Value = 11.7340452250213
$ws.Range("C6").Value = 8.139747843071419
$ws.Range("E6").Value = 11.25378620748242
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 28.16659315773495
$ws.Range("H6").Value = 13.76793503523901
$ws.Range("M6").Value = 14.09768114802013
$ws.Range("B7").Value = 11.9215798285023
$ws.Range("C7").Value = 8.288312776528485
$ws.Range("E7").Value = 11.27906653962753
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 28.23931776140495
$ws.Range("H7").Value = 13.75285418608073
$ws.Range("M7").Value = 14.1801199467284
$ws.Range("B8").Value = 12.72007924577467
$ws.Range("C8").Value = 8.910211494271607
$ws.Range("E8").Value = 11.39908182115308
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 28.59712543631014
$ws.Range("H8").Value = 13.69615991251089
$ws.Range("M8").Value = 14.54680761348506
$ws.Range("B9").Value = 14.1684931557617
$ws.Range("C9").Value = 10.00504629876744
$ws.Range("E9").Value = 11.66341149670066
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 29.42611008283443
$ws.Range("H9").Value = 13.62094219649993
$ws.Range("M9").Value = 15.26955401593806
$ws.Range("B10").Value = 15.1520426716192
$ws.Range("C10").Value = 10.73116043948559
$ws.Range("E10").Value = 11.87382801656629
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 30.10692747476487
$ws.Range("H10").Value = 13.58806643630854
$ws.Range("M10").Value = 15.79697538106704
$ws.Range("B11").Value = 15.58055831164713
$ws.Range("C11").Value = 11.04419003124577
$ws.Range("E11").Value = 11.97272910688389
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 30.43077750536516
$ws.Range("H11").Value = 13.57806199951176
$ws.Range("M11").Value = 16.03497517667888
$ws.Range("B12").Value = 15.74001447678479
$ws.Range("C12").Value = 11.16022330005124
$ws.Range("E12").Value = 12.01060721764728
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 30.55531268911785
$ws.Range("H12").Value = 13.57499217871582
$ws.Range("M12").Value = 16.12473178644266
$ws.Range("B13").Value = 15.70579920070314
$ws.Range("C13").Value = 11.13534508318827
$ws.Range("E13").Value = 12.00243102827772
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 30.52840947622428
$ws.Range("H13").Value = 13.57562125648984
$ws.Range("M13").Value = 16.10541890003152
$ws.Range("B14").Value = 15.59373373559851
$ws.Range("C14").Value = 11.05378643692784
$ws.Range("E14").Value = 11.97583699711597
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 30.44098578977034
$ws.Range("H14").Value = 13.57779500428979
$ws.Range("M14").Value = 16.04236740418987
$ws.Range("B15").Value = 15.52472148323848
$ws.Range("C15").Value = 11.00350281173471
$ws.Range("E15").Value = 11.95960197877825
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 30.38767966973788
$ws.Range("H15").Value = 13.57922026617152
$ws.Range("M15").Value = 16.00369588673704
$ws.Range("B16").Value = 15.12364924704883
$ws.Range("C16").Value = 10.71035369591979
$ws.Range("E16").Value = 11.86742598262144
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 30.08603641341809
$ws.Range("H16").Value = 13.58882045994952
$ws.Range("M16").Value = 15.78137481044645
$ws.Range("B17").Value = 14.87268931712212
$ws.Range("C17").Value = 10.52607435096142
$ws.Range("E17").Value = 11.81167205604667
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 29.90451714244323
$ws.Range("H17").Value = 13.59598285057557
$ws.Range("M17").Value = 15.64442771034449
$ws.Range("B18").Value = 14.72656878988027
$ws.Range("C18").Value = 10.41845571174456
$ws.Range("E18").Value = 11.77990526567071
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 29.80145164371953
$ws.Range("H18").Value = 13.60056810238132
$ws.Range("M18").Value = 15.56548317306933
$ws.Range("B19").Value = 14.67679319419776
$ws.Range("C19").Value = 10.38173928072447
$ws.Range("E19").Value = 11.7692022725292
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 29.7667894891748
$ws.Range("H19").Value = 13.60220035008073
$ws.Range("M19").Value = 15.53872652440257
$ws.Range("B20").Value = 14.89958887647069
$ws.Range("C20").Value = 10.545859643639
$ws.Range("E20").Value = 11.81757618589228
$ws.Range("F20").Value = 20.24955283636157
$ws.Range("G20").Value = 29.92370248720716
$ws.Range("H20").Value = 13.59517216045596
$ws.Range("M20").Value = 15.65902487780933
$ws.Range("B21").Value = 15.62672710804135
$ws.Range("C21").Value = 11.07781027802546
$ws.Range("E21").Value = 11.98363698813959
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 30.46661376388956
$ws.Range("H21").Value = 13.57713696725541
$ws.Range("M21").Value = 16.06089786990551
$ws.Range("B22").Value = 16.08552244286131
$ws.Range("C22").Value = 11.41086880334843
$ws.Range("E22").Value = 12.09463779929804
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 30.83244836623184
$ws.Range("H22").Value = 13.56954146043679
$ws.Range("M22").Value = 16.32135689145433
$ws.Range("B23").Value = 15.84218538816967
$ws.Range("C23").Value = 11.23445023715861
$ws.Range("E23").Value = 12.03517901666438
$ws.Range("F23").Value = 21.82633154475857
$ws.Range("G23").Value = 30.63623341580077
$ws.Range("H23").Value = 13.57320971993635
$ws.Range("M23").Value = 16.18257418656928
$ws.Range("B24").Value = 14.88743331052833
$ws.Range("C24").Value = 10.53691992552173
$ws.Range("E24").Value = 11.81490603285971
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 29.91502475726653
$ws.Range("H24").Value = 13.59553721772575
$ws.Range("M24").Value = 15.6524261511686
$ws.Range("B25").Value = 13.79025907605309
$ws.Range("C25").Value = 9.722480706657988
$ws.Range("E25").Value = 11.58894395036153
$ws.Range("F25").Value = 18.34778573295691
$ws.Range("G25").Value = 29.18877964116062
$ws.Range("H25").Value = 13.63739121096783
$ws.Range("M25").Value = 15.07425263176856
